$p = $ppt.ActivePresentation

# 1. Replace the "expectation" bullet on slide 2 with the research question.
$s2 = $p.Slides.Item(2)
$contentPh = $s2.Shapes.Item(2)
$para = $contentPh.TextFrame.TextRange.Paragraphs(4, 1)
$para.Text = "What is the effect of COVID-19 on tweet sentiment?"

# 2. Reorder the closing slides: move the "Thank you" slide ahead of "References".
$p.Slides.Item(6).MoveTo(5)

# 3. Refresh the cached slide-number fields on the two slides that moved so the
#    printed numbers match their new positions (5 and 6).
$s5 = $p.Slides.Item(5)
$s6 = $p.Slides.Item(6)
$s5.HeadersFooters.SlideNumber.Visible = $false
$s5.HeadersFooters.SlideNumber.Visible = $true
$s6.HeadersFooters.SlideNumber.Visible = $false
$s6.HeadersFooters.SlideNumber.Visible = $true
